# Venezuela Primera Division - base update 03-06-2024
# The underlying match records (columns B:AD) were reshuffled between
# several rows while each row kept its own rank/position (column A).
# Capture the "before" B:AD values of every involved row first, then
# write them back in their new positions, so the cyclic swaps don't
# clobber data we still need to read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return $ws.Range("B$row`:AD$row").Value2
}

function Set-RowData($row, $data) {
    $ws.Range("B$row`:AD$row").Value2 = $data
}

# Snapshot the original data for every row that participates in a swap.
$rows = 93, 94, 95, 96, 98, 99, 102, 103, 116, 117, 135, 136, 157, 158
$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = Get-RowData $r
}

# Mapping: new content of row <key> = old content of row <value>
$srcMap = @{
    93  = 96
    94  = 93
    96  = 94
    95  = 99
    99  = 98
    98  = 95
    102 = 103
    103 = 102
    116 = 117
    117 = 116
    135 = 136
    136 = 135
    157 = 158
    158 = 157
}

foreach ($r in $rows) {
    Set-RowData $r $orig[$srcMap[$r]]
}
